# Fix spelling error "Architectur" -> "Architektur" in the
# "Model view client Architectur -> GUI konstruieren." bullet on the
# "Technology" slide (Content Placeholder 2).
#
# In the original file this phrase is split across three runs:
#   " " | "Architectur" | " -> GUI konstruieren. "
# The fix merges them into a single corrected run:
#   " Architektur -> GUI konstruieren. "

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

$found = $tr.Find("Architectur -> GUI konstruieren. ")

$start = $found.Start
$len = $found.Length

# Include the leading space run that precedes "Architectur" so the whole
# segment (leading space + misspelled word + trailing text) is replaced
# by one corrected run, matching the target XML structure.
$segment = $tr.Characters($start - 1, $len + 1)
$segment.Text = " Architektur -> GUI konstruieren. "
